$d = $word.ActiveDocument

# Merge the split "ТЕЛ" label runs into a single run (does not cross the
# following <w:br/>, so it stays its own run)
$d.Content.Find.Execute("Т" + "ЕЛ" + ".: +7 (952) 548-88-10", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ТЕЛ.: +7 (952) 548-88-10", 2)

# Merge the split "ПОЧТА" label runs into a single run; the leading <w:br/>
# (in its own/previous run) is left untouched since the match text starts
# right after it
$d.Content.Find.Execute("П" + "ОЧТА" + ": mixailtrifonov@yandex.ru", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ПОЧТА: mixailtrifonov@yandex.ru", 2)

# Merge the split "ПУБЛИКАЦИИ" label runs into a single run
$d.Content.Find.Execute("П" + "УБЛИКАЦИИ" + ":", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ПУБЛИКАЦИИ:", 2)

# Remove the trailing placeholder text, merging the "ВРЕМЯ..." run with the
# former ": день, дата, ауд." run and dropping the placeholder
$d.Content.Find.Execute("ВРЕМЯ РАБОТЫ СО СТУДЕНТАМИ" + ": день, дата, ауд.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ВРЕМЯ РАБОТЫ СО СТУДЕНТАМИ:", 2)
